$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Fecha" (column D) values for rows 20-107, and a handful of
# "Origen" (column O) corrections, per the source diff. ---

$ws.Range("D20").Value = 44425
$ws.Range("D21").Value = 44425
$ws.Range("D22").Value = 44252
$ws.Range("D23").Value = 44252
$ws.Range("D24").Value = 44322
$ws.Range("D25").Value = 44322
$ws.Range("D26").Value = 44166
$ws.Range("D27").Value = 44166
$ws.Range("D28").Value = 44308
$ws.Range("D29").Value = 44308
$ws.Range("D30").Value = 44239
$ws.Range("D31").Value = 44239
$ws.Range("D32").Value = 44306
$ws.Range("D33").Value = 44306
$ws.Range("D34").Value = 44271
$ws.Range("D35").Value = 44271
$ws.Range("D36").Value = 44280
$ws.Range("D37").Value = 44280
$ws.Range("D38").Value = 44299
$ws.Range("D39").Value = 44299
$ws.Range("D40").Value = 44362
$ws.Range("D41").Value = 44362
$ws.Range("D42").Value = 44204
$ws.Range("D43").Value = 44204
$ws.Range("D44").Value = 44427
$ws.Range("D45").Value = 44427
$ws.Range("D46").Value = 44222
$ws.Range("D47").Value = 44222
$ws.Range("D48").Value = 44237
$ws.Range("D49").Value = 44237
$ws.Range("D50").Value = 44257
$ws.Range("D51").Value = 44257
$ws.Range("D52").Value = 44194
$ws.Range("D53").Value = 44194
$ws.Range("D54").Value = 44383
$ws.Range("D55").Value = 44383
$ws.Range("D56").Value = 44169
$ws.Range("D57").Value = 44169
$ws.Range("D58").Value = 44336
$ws.Range("D59").Value = 44336
$ws.Range("D60").Value = 44371
$ws.Range("O60").Value = 'Región de Ñuble'
$ws.Range("D61").Value = 44371
$ws.Range("O61").Value = 'Región de Ñuble'
$ws.Range("D62").Value = 44274
$ws.Range("D63").Value = 44274
$ws.Range("D64").Value = 44320
$ws.Range("O64").Value = 'Región Metropolitana'
$ws.Range("D65").Value = 44320
$ws.Range("O65").Value = 'Región Metropolitana'
$ws.Range("D66").Value = 44405
$ws.Range("D67").Value = 44405
$ws.Range("D68").Value = 44224
$ws.Range("D69").Value = 44224
$ws.Range("D70").Value = 44327
$ws.Range("D71").Value = 44327
$ws.Range("D72").Value = 44209
$ws.Range("D73").Value = 44209
$ws.Range("D74").Value = 44231
$ws.Range("D75").Value = 44231
$ws.Range("D76").Value = 44313
$ws.Range("D77").Value = 44313
$ws.Range("D78").Value = 44330
$ws.Range("D79").Value = 44330
$ws.Range("D80").Value = 44391
$ws.Range("D81").Value = 44391
$ws.Range("D82").Value = 44350
$ws.Range("D83").Value = 44350
$ws.Range("D84").Value = 44278
$ws.Range("D85").Value = 44278
$ws.Range("D86").Value = 44358
$ws.Range("D87").Value = 44358
$ws.Range("D88").Value = 44250
$ws.Range("D89").Value = 44250
$ws.Range("D90").Value = 44292
$ws.Range("D91").Value = 44292
$ws.Range("D92").Value = 44420
$ws.Range("D93").Value = 44420
$ws.Range("D94").Value = 44245
$ws.Range("D95").Value = 44245
$ws.Range("D96").Value = 44161
$ws.Range("O96").Value = 'Región de Ñuble'
$ws.Range("D97").Value = 44161
$ws.Range("O97").Value = 'Región de Ñuble'
$ws.Range("D98").Value = 44159
$ws.Range("D99").Value = 44159
$ws.Range("D100").Value = 44344
$ws.Range("D101").Value = 44344
$ws.Range("D102").Value = 44316
$ws.Range("O102").Value = 'Región Metropolitana'
$ws.Range("D103").Value = 44316
$ws.Range("O103").Value = 'Región Metropolitana'
$ws.Range("D104").Value = 44398
$ws.Range("D105").Value = 44398
$ws.Range("D106").Value = 44217
$ws.Range("D107").Value = 44217

# --- Append 6 new rows (108-113) with the same record layout used
# throughout the sheet. ---

# Row 108
$ws.Range("A108").Value = 11
$ws.Range("B108").Value = 'Vega Monumental Concepción'
$ws.Range("C108").Value = 'Bíobío'
$ws.Range("D108").Value = 44376
$ws.Range("E108").Value = 8
$ws.Range("F108").Value = 100112040
$ws.Range("G108").Value = 'Cilantro'
$ws.Range("H108").Value = 'Sin especificar'
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 200
$ws.Range("K108").Value = 600
$ws.Range("L108").Value = 700
$ws.Range("M108").Value = 650
$ws.Range("N108").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O108").Value = 'Región de Ñuble'
$ws.Range("P108").Value = 650
$ws.Range("Q108").Value = 1
$ws.Range("R108").Value = 'Hortaliza'

# Row 109
$ws.Range("A109").Value = 11
$ws.Range("B109").Value = 'Vega Monumental Concepción'
$ws.Range("C109").Value = 'Bíobío'
$ws.Range("D109").Value = 44376
$ws.Range("E109").Value = 8
$ws.Range("F109").Value = 100112040
$ws.Range("G109").Value = 'Cilantro'
$ws.Range("H109").Value = 'Sin especificar'
$ws.Range("I109").Value = 'Segunda'
$ws.Range("J109").Value = 100
$ws.Range("K109").Value = 500
$ws.Range("L109").Value = 500
$ws.Range("M109").Value = 500
$ws.Range("N109").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O109").Value = 'Región de Ñuble'
$ws.Range("P109").Value = 500
$ws.Range("Q109").Value = 1
$ws.Range("R109").Value = 'Hortaliza'

# Row 110
$ws.Range("A110").Value = 11
$ws.Range("B110").Value = 'Vega Monumental Concepción'
$ws.Range("C110").Value = 'Bíobío'
$ws.Range("D110").Value = 44334
$ws.Range("E110").Value = 8
$ws.Range("F110").Value = 100112040
$ws.Range("G110").Value = 'Cilantro'
$ws.Range("H110").Value = 'Sin especificar'
$ws.Range("I110").Value = 'Primera'
$ws.Range("J110").Value = 200
$ws.Range("K110").Value = 600
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 650
$ws.Range("N110").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O110").Value = 'Región de Ñuble'
$ws.Range("P110").Value = 650
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = 'Hortaliza'

# Row 111
$ws.Range("A111").Value = 11
$ws.Range("B111").Value = 'Vega Monumental Concepción'
$ws.Range("C111").Value = 'Bíobío'
$ws.Range("D111").Value = 44334
$ws.Range("E111").Value = 8
$ws.Range("F111").Value = 100112040
$ws.Range("G111").Value = 'Cilantro'
$ws.Range("H111").Value = 'Sin especificar'
$ws.Range("I111").Value = 'Segunda'
$ws.Range("J111").Value = 100
$ws.Range("K111").Value = 500
$ws.Range("L111").Value = 500
$ws.Range("M111").Value = 500
$ws.Range("N111").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O111").Value = 'Región de Ñuble'
$ws.Range("P111").Value = 500
$ws.Range("Q111").Value = 1
$ws.Range("R111").Value = 'Hortaliza'

# Row 112
$ws.Range("A112").Value = 11
$ws.Range("B112").Value = 'Vega Monumental Concepción'
$ws.Range("C112").Value = 'Bíobío'
$ws.Range("D112").Value = 44168
$ws.Range("E112").Value = 8
$ws.Range("F112").Value = 100112040
$ws.Range("G112").Value = 'Cilantro'
$ws.Range("H112").Value = 'Sin especificar'
$ws.Range("I112").Value = 'Primera'
$ws.Range("J112").Value = 200
$ws.Range("K112").Value = 600
$ws.Range("L112").Value = 700
$ws.Range("M112").Value = 650
$ws.Range("N112").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O112").Value = 'Región de Ñuble'
$ws.Range("P112").Value = 650
$ws.Range("Q112").Value = 1
$ws.Range("R112").Value = 'Hortaliza'

# Row 113
$ws.Range("A113").Value = 11
$ws.Range("B113").Value = 'Vega Monumental Concepción'
$ws.Range("C113").Value = 'Bíobío'
$ws.Range("D113").Value = 44168
$ws.Range("E113").Value = 8
$ws.Range("F113").Value = 100112040
$ws.Range("G113").Value = 'Cilantro'
$ws.Range("H113").Value = 'Sin especificar'
$ws.Range("I113").Value = 'Segunda'
$ws.Range("J113").Value = 100
$ws.Range("K113").Value = 500
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = 500
$ws.Range("N113").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O113").Value = 'Región de Ñuble'
$ws.Range("P113").Value = 500
$ws.Range("Q113").Value = 1
$ws.Range("R113").Value = 'Hortaliza'

# Ensure the newly appended "Fecha" cells use the same date number format
# as the rest of column D.
$ws.Range("D108:D113").NumberFormat = "YYYY-MM-DD HH:MM:SS"

